# draft-gandhi-spring-sr-enhanced-plm-04.pptx -- "Add files via upload"
#
# Three text tweaks in this revision:
#  1. Handout master "date" field cache: 2/9/21 -> 2/10/21
#  2. Slide 10, content placeholder: "Timestamp label (TBA1) is defined ..."
#     -> "Timestamp labels (TBA1 and TBA2) are defined ..."
#  3. Slide 11, content placeholder: "Timestamp Endpoint Function " ->
#     "Timestamp Endpoint Functions " and " (TBA3) is defined ..." ->
#     " (TBA3 and TBA4) is defined ..."

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Handout master date placeholder (cached text of the datetimeFigureOut
#    field). Best-effort: some hosts treat the handout/notes master as a
#    read-only part, so guard this so the rest of the edits still apply.
# ---------------------------------------------------------------------
try {
    $hm = $p.HandoutMaster
    $dateShape = $hm.Shapes.Item(2)
    $dateRange = $dateShape.TextFrame.TextRange
    $oldDate = "2/9/21"
    $newDate = "2/10/21"
    if ($dateRange.Text -eq $oldDate) {
        $dateChars = $dateRange.Characters(1, $oldDate.Length)
        $dateChars.Text = $newDate
    }
} catch {
    Write-Output "handout master date field not editable: $_"
}

# ---------------------------------------------------------------------
# 2) Slide 10 - "Timestamp label (TBA1) is defined for Timestamp, Pop and
#    Forward function"
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(5)
$tr10 = $shp10.TextFrame.TextRange

$old10 = "Timestamp label (TBA1) is defined for Timestamp, Pop and Forward function"
$new10 = "Timestamp labels (TBA1 and TBA2) are defined for Timestamp, Pop and Forward function"

$full10 = $tr10.Text
$idx10 = $full10.IndexOf($old10)
if ($idx10 -ge 0) {
    $run10 = $tr10.Characters($idx10 + 1, $old10.Length)
    $run10.Text = $new10
}

# ---------------------------------------------------------------------
# 3) Slide 11 - "Timestamp Endpoint Function End.TSF (TBA3) is defined for
#    Timestamp and Forward and is carried with the Session-Reflector node
#    SID"
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(5)
$tr11 = $shp11.TextFrame.TextRange

$oldA = " (TBA3) is defined for Timestamp and Forward and is carried with the Session-Reflector node SID"
$newA = " (TBA3 and TBA4) is defined for Timestamp and Forward and is carried with the Session-Reflector node SID"

$oldB = "Timestamp Endpoint Function "
$newB = "Timestamp Endpoint Functions "

# Replace the later occurrence first so the earlier index is unaffected by
# any length change.
$full11 = $tr11.Text
$idxA = $full11.IndexOf($oldA)
if ($idxA -ge 0) {
    $runA = $tr11.Characters($idxA + 1, $oldA.Length)
    $runA.Text = $newA
}

$full11 = $tr11.Text
$idxB = $full11.IndexOf($oldB)
if ($idxB -ge 0) {
    $runB = $tr11.Characters($idxB + 1, $oldB.Length)
    $runB.Text = $newB
}
